$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2648451569472172
$ws.Range("C2").Value = 0.06939515048310341
$ws.Range("E2").Value = 0.5631802312988512
$ws.Range("F2").Value = 2.073558812778955
$ws.Range("G2").Value = 0.4535287523932894
$ws.Range("H2").Value = 0.621328140970526
$ws.Range("J2").Value = 0.04456917779960357
$ws.Range("K2").Value = 0.2577744102338499
$ws.Range("M2").Value = 0.4017623662331715
$ws.Range("O2").Value = 2.105377801788336

$ws.Range("B3").Value = 0.2321411886632347
$ws.Range("C3").Value = 0.06688639323184731
$ws.Range("E3").Value = 0.5528305900229853
$ws.Range("F3").Value = 2.065961215728322
$ws.Range("G3").Value = 0.4590390420931456
$ws.Range("H3").Value = 0.6274776546296081
$ws.Range("J3").Value = 0.04471664244262641
$ws.Range("K3").Value = 0.2252679536051119
$ws.Range("M3").Value = 0.3786317057845139
$ws.Range("O3").Value = 2.129708494002685

$ws.Range("B4").Value = 0.212012630724729
$ws.Range("C4").Value = 0.06534668349953421
$ws.Range("E4").Value = 0.5467358778822273
$ws.Range("F4").Value = 2.062461071519834
$ws.Range("G4").Value = 0.4627584387925623
$ws.Range("H4").Value = 0.6315253166097463
$ws.Range("J4").Value = 0.04481790224832416
$ws.Range("K4").Value = 0.2052215070626033
$ws.Range("M4").Value = 0.3645554958694888
$ws.Range("O4").Value = 2.145921564995248

$ws.Range("B5").Value = 0.2037984957325136
$ws.Range("C5").Value = 0.06471947030497205
$ws.Range("E5").Value = 0.5443177112772162
$ws.Range("F5").Value = 2.061327672639166
$ws.Range("G5").Value = 0.4643585217900252
$ws.Range("H5").Value = 0.6332431509224179
$ws.Range("J5").Value = 0.04486186884625099
$ws.Range("K5").Value = 0.1970309852239041
$ws.Range("M5").Value = 0.3588513314354103
$ws.Range("O5").Value = 2.152848683791177

$ws.Range("B6").Value = 0.2024338609564325
$ws.Range("C6").Value = 0.06461533777898154
$ws.Range("E6").Value = 0.5439201341323709
$ws.Range("F6").Value = 2.06115716650919
$ws.Range("G6").Value = 0.4646293087253497
$ws.Range("H6").Value = 0.6335325266417513
$ws.Range("J6").Value = 0.04486933291460993
$ws.Range("K6").Value = 0.1956696765995787
$ws.Range("M6").Value = 0.3579061010055611
$ws.Range("O6").Value = 2.154018258272032

$ws.Range("B7").Value = 0.2119018982236867
$ws.Range("C7").Value = 0.06533822366652942
$ws.Range("E7").Value = 0.5467030003954392
$ws.Range("F7").Value = 2.062444599962888
$ws.Range("G7").Value = 0.462779676485134
$ws.Range("H7").Value = 0.6315482070281107
$ws.Range("J7").Value = 0.04481848424526547
$ws.Range("K7").Value = 0.2051111328900674
$ws.Range("M7").Value = 0.3644784375338901
$ws.Range("O7").Value = 2.146013690461942

$ws.Range("B8").Value = 0.2535791398760523
$ws.Range("C8").Value = 0.06853002588876222
$ws.Range("E8").Value = 0.5595577634227169
$ws.Range("F8").Value = 2.070697418761767
$ws.Range("G8").Value = 0.4553588885072557
$ws.Range("H8").Value = 0.6233920884511903
$ws.Range("J8").Value = 0.04461780417199179
$ws.Range("K8").Value = 0.2465846230229261
$ws.Range("M8").Value = 0.3937609092317516
$ws.Range("O8").Value = 2.113502514398675

$ws.Range("B9").Value = 0.3349062129612435
$ws.Range("C9").Value = 0.07479216954615708
$ws.Range("E9").Value = 0.5868266450252548
$ws.Range("F9").Value = 2.096125596680082
$ws.Range("G9").Value = 0.4434779939266846
$ws.Range("H9").Value = 0.6095537205753701
$ws.Range("J9").Value = 0.04430896245136928
$ws.Range("K9").Value = 0.3272013342665332
$ws.Range("M9").Value = 0.4521745746135437
$ws.Range("O9").Value = 2.059864243709299

$ws.Range("B10").Value = 0.3943917719728063
$ws.Range("C10").Value = 0.07939213296144487
$ws.Range("E10").Value = 0.6081170759557324
$ws.Range("F10").Value = 2.12045188888753
$ws.Range("G10").Value = 0.4363842403982616
$ws.Range("H10").Value = 0.6006992173518242
$ws.Range("J10").Value = 0.04413324813570618
$ws.Range("K10").Value = 0.3859752718024936
$ws.Range("M10").Value = 0.4956865967892767
$ws.Range("O10").Value = 2.026634219022625

$ws.Range("B11").Value = 0.4213916910273952
$ws.Range("C11").Value = 0.08148406107864048
$ws.Range("E11").Value = 0.6180754916050546
$ws.Range("F11").Value = 2.132746334632387
$ws.Range("G11").Value = 0.4335136019092616
$ws.Range("H11").Value = 0.5969557351015524
$ws.Range("J11").Value = 0.04406433556560962
$ws.Range("K11").Value = 0.412610207403219
$ws.Range("M11").Value = 0.5156091276955124
$ws.Range("O11").Value = 2.012860911190032

$ws.Range("B12").Value = 0.4316067174264901
$ws.Range("C12").Value = 0.0822760768187436
$ws.Range("E12").Value = 0.6218857319189937
$ws.Range("F12").Value = 2.137578617892672
$ws.Range("G12").Value = 0.4324779267023615
$ws.Range("H12").Value = 0.5955790630197342
$ws.Range("J12").Value = 0.04403981739016238
$ws.Range("K12").Value = 0.4226810903254261
$ws.Range("M12").Value = 0.5231715451259475
$ws.Range("O12").Value = 2.007838683217471

$ws.Range("B13").Value = 0.4294071495960168
$ws.Range("C13").Value = 0.08210550981951314
$ws.Range("E13").Value = 0.6210633863714463
$ws.Range("F13").Value = 2.136530042714895
$ws.Range("G13").Value = 0.4326986909657577
$ws.Range("H13").Value = 0.5958737348529084
$ws.Range("J13").Value = 0.0440450277725084
$ws.Range("K13").Value = 0.4205128301776426
$ws.Range("M13").Value = 0.521542039406782
$ws.Range("O13").Value = 2.008911702887559

$ws.Range("B14").Value = 0.4222322762189492
$ws.Range("C14").Value = 0.08154922409721621
$ws.Range("E14").Value = 0.6183881771952997
$ws.Range("F14").Value = 2.133140349195685
$ws.Range("G14").Value = 0.433427365940048
$ws.Range("H14").Value = 0.5968416557173128
$ws.Range("J14").Value = 0.04406228686625724
$ws.Range("K14").Value = 0.4134390527040068
$ws.Range("M14").Value = 0.5162309296515843
$ws.Range("O14").Value = 2.012443851381434

$ws.Range("B15").Value = 0.4178362364794168
$ws.Range("C15").Value = 0.08120846136456805
$ws.Range("E15").Value = 0.6167546370973298
$ws.Range("F15").Value = 2.131087068704488
$ws.Range("G15").Value = 0.4338803946039178
$ws.Range("H15").Value = 0.5974398619146868
$ws.Range("J15").Value = 0.04407306378451104
$ws.Range("K15").Value = 0.4091041646325948
$ws.Range("M15").Value = 0.5129800806488021
$ws.Range("O15").Value = 2.014632592075358

$ws.Range("B16").Value = 0.3926259997705586
$ws.Range("C16").Value = 0.07925540269033604
$ws.Range("E16").Value = 0.6074717613158782
$ws.Range("F16").Value = 2.119673135770938
$ws.Range("G16").Value = 0.4365790251919535
$ws.Range("H16").Value = 0.6009495867633063
$ws.Range("J16").Value = 0.04413797287896593
$ws.Range("K16").Value = 0.384232521915294
$ws.Range("M16").Value = 0.4943871749426805
$ws.Range("O16").Value = 2.027561389141894

$ws.Range("B17").Value = 0.3771444650262481
$ws.Range("C17").Value = 0.07805706101929388
$ws.Range("E17").Value = 0.6018469469708521
$ws.Range("F17").Value = 2.112985664998931
$ws.Range("G17").Value = 0.4383259098841776
$ws.Range("H17").Value = 0.6031755474983385
$ws.Range("J17").Value = 0.04418061000967022
$ws.Range("K17").Value = 0.3689481359554065
$ws.Range("M17").Value = 0.4830137642598729
$ws.Range("O17").Value = 2.035837008103613

$ws.Range("B18").Value = 0.3682342487029757
$ws.Range("C18").Value = 0.07736775213091107
$ws.Range("E18").Value = 0.5986374317034517
$ws.Range("F18").Value = 2.10925482055309
$ws.Range("G18").Value = 0.4393642048834039
$ws.Range("H18").Value = 0.6044826376472017
$ws.Range("J18").Value = 0.04420617151534678
$ws.Range("K18").Value = 0.3601474303827104
$ws.Range("M18").Value = 0.4764842115490566
$ws.Range("O18").Value = 2.040723336942975

$ws.Range("B19").Value = 0.3652164503007782
$ws.Range("C19").Value = 0.07713435656314971
$ws.Range("E19").Value = 0.5975551668750541
$ws.Range("F19").Value = 2.108011476666533
$ws.Range("G19").Value = 0.4397215090677093
$ws.Range("H19").Value = 0.6049297949100136
$ws.Range("J19").Value = 0.0442150046590104
$ws.Range("K19").Value = 0.3571660433196939
$ws.Range("M19").Value = 0.4742755106425847
$ws.Range("O19").Value = 2.042399466879701

$ws.Range("B20").Value = 0.3787930894614817
$ws.Range("C20").Value = 0.07818463263097897
$ws.Range("E20").Value = 0.6024430558725555
$ws.Range("F20").Value = 2.113685592160948
$ws.Range("G20").Value = 0.438136479520324
$ws.Range("H20").Value = 0.6029358189482181
$ws.Range("J20").Value = 0.04417596386318223
$ws.Range("K20").Value = 0.3705761755383321
$ws.Range("M20").Value = 0.4842232299051901
$ws.Range("O20").Value = 2.034942969206554

$ws.Range("B21").Value = 0.4243399662729246
$ws.Range("C21").Value = 0.08171262332518836
$ws.Range("E21").Value = 0.6191728877349192
$ws.Range("F21").Value = 2.134131190170663
$ws.Range("G21").Value = 0.4332119411379978
$ws.Range("H21").Value = 0.5965562438196343
$ws.Range("J21").Value = 0.04405717469620285
$ws.Range("K21").Value = 0.4155172094236264
$ws.Range("M21").Value = 0.5177904402685556
$ws.Range("O21").Value = 2.011401122070666

$ws.Range("B22").Value = 0.4540532601330654
$ws.Range("C22").Value = 0.08401745405038241
$ws.Range("E22").Value = 0.6303352695513098
$ws.Range("F22").Value = 2.148523190698398
$ws.Range("G22").Value = 0.4302929738868002
$ws.Range("H22").Value = 0.5926252374063878
$ws.Range("J22").Value = 0.04398873076637422
$ws.Range("K22").Value = 0.4447999737770942
$ws.Range("M22").Value = 0.5398344148196514
$ws.Range("O22").Value = 1.997142718542278

$ws.Range("B23").Value = 0.4381998676004457
$ws.Range("C23").Value = 0.082787426981767
$ws.Range("E23").Value = 0.6243568225631577
$ws.Range("F23").Value = 2.140747693674001
$ws.Range("G23").Value = 0.431823433891573
$ws.Range("H23").Value = 0.5947014757540998
$ws.Range("J23").Value = 0.04402442192577105
$ws.Range("K23").Value = 0.4291795260426738
$ws.Range("M23").Value = 0.5280595495571418
$ws.Range("O23").Value = 2.004649435022415

$ws.Range("B24").Value = 0.3780477764501882
$ws.Range("C24").Value = 0.07812695864262764
$ws.Range("E24").Value = 0.6021734794461366
$ws.Range("F24").Value = 2.113368800388457
$ws.Range("G24").Value = 0.4382220151282894
$ws.Range("H24").Value = 0.6030441150133896
$ws.Range("J24").Value = 0.0441780611183944
$ws.Range("K24").Value = 0.3698401808774179
$ws.Range("M24").Value = 0.4836764018920121
$ws.Range("O24").Value = 2.035346763738332

$ws.Range("B25").Value = 0.3129502306972256
$ws.Range("C25").Value = 0.07309804842563494
$ws.Range("E25").Value = 0.5792291053882579
$ws.Range("F25").Value = 2.088256034321262
$ws.Range("G25").Value = 0.4464054232406411
$ws.Range("H25").Value = 0.6130667501616074
$ws.Range("J25").Value = 0.04438349104491124
$ws.Range("K25").Value = 0.3054708247636029
$ws.Range("M25").Value = 0.4362669390882417
$ws.Range("O25").Value = 2.073290685977142
